$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.516.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.959.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4798"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4081"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.99"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08501"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.060"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.49"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.957.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.597"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.173"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.012"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.65"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001074"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06623"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.851"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.504.85"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.175.67"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.36"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.178"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.842"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.55"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9882"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09675"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.645"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.694"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.166"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02338"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06197"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.257"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6250"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.22"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.009"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1922"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.348"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5964"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.06%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.413"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06817"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.25"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.06%  "
